# This script applies the "Added easy replay battle testing" edit:
# appends 34 new card rows (rows 64-97) to Sheet1 -- the Life splinter
# cards first, then the Death + Dragon splinter cards -- entered a
# column at a time (Card, Role, Element, AttackType, Ability1,
# Ability2) the way the author actually typed/pasted them, which is
# what the shared-string ordering in the saved file reflects. Also
# nudges a couple of workbook/window display settings to match the
# author's working session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Life splinter cards (rows 64-79) ---

# Card
$ws.Range("A64").Value = 'Divine Healer'
$ws.Range("A65").Value = 'Feral Spirit'
$ws.Range("A66").Value = 'Silvershield Knight'
$ws.Range("A67").Value = 'Silvershield Warrior'
$ws.Range("A68").Value = 'Cave Slug'
$ws.Range("A69").Value = 'Crystal Jaguar'
$ws.Range("A70").Value = 'Lone Boatman'
$ws.Range("A71").Value = 'Herbalist'
$ws.Range("A72").Value = 'Tyrus Paladium'
$ws.Range("A73").Value = 'Peacebringer'
$ws.Range("A74").Value = 'Silvershield Paladin'
$ws.Range("A75").Value = 'Clay Golem'
$ws.Range("A76").Value = 'Truthspeaker'
$ws.Range("A77").Value = 'Luminous Eagle'
$ws.Range("A78").Value = 'Shieldbearer'
$ws.Range("A79").Value = 'Mother Khala'

# Role
$ws.Range("B64").Value = 'monster'
$ws.Range("B65").Value = 'monster'
$ws.Range("B66").Value = 'monster'
$ws.Range("B67").Value = 'monster'
$ws.Range("B68").Value = 'monster'
$ws.Range("B69").Value = 'monster'
$ws.Range("B70").Value = 'monster'
$ws.Range("B71").Value = 'monster'
$ws.Range("B72").Value = 'summoner'
$ws.Range("B73").Value = 'monster'
$ws.Range("B74").Value = 'monster'
$ws.Range("B75").Value = 'monster'
$ws.Range("B76").Value = 'monster'
$ws.Range("B77").Value = 'monster'
$ws.Range("B78").Value = 'monster'
$ws.Range("B79").Value = 'summoner'

# Element
$ws.Range("C64").Value = 'life'
$ws.Range("C65").Value = 'life'
$ws.Range("C66").Value = 'life'
$ws.Range("C67").Value = 'life'
$ws.Range("C68").Value = 'life'
$ws.Range("C69").Value = 'life'
$ws.Range("C70").Value = 'life'
$ws.Range("C71").Value = 'life'
$ws.Range("C72").Value = 'life'
$ws.Range("C73").Value = 'life'
$ws.Range("C74").Value = 'life'
$ws.Range("C75").Value = 'life'
$ws.Range("C76").Value = 'life'
$ws.Range("C77").Value = 'life'
$ws.Range("C78").Value = 'life'
$ws.Range("C79").Value = 'life'

# ManaCost
$ws.Range("D64").Value = 3
$ws.Range("D65").Value = 3
$ws.Range("D66").Value = 6
$ws.Range("D67").Value = 4
$ws.Range("D68").Value = 5
$ws.Range("D69").Value = 4
$ws.Range("D70").Value = 5
$ws.Range("D71").Value = 2
$ws.Range("D72").Value = 3
$ws.Range("D73").Value = 4
$ws.Range("D74").Value = 5
$ws.Range("D75").Value = 6
$ws.Range("D76").Value = 3
$ws.Range("D77").Value = 6
$ws.Range("D78").Value = 8
$ws.Range("D79").Value = 3

# Dmg
$ws.Range("E64").Value = 0
$ws.Range("E65").Value = 1
$ws.Range("E66").Value = 1
$ws.Range("E67").Value = 1
$ws.Range("E68").Value = 2
$ws.Range("E69").Value = 1
$ws.Range("E70").Value = 2
$ws.Range("E71").Value = 1
$ws.Range("E72").Value = 0
$ws.Range("E73").Value = 2
$ws.Range("E74").Value = 1
$ws.Range("E75").Value = 3
$ws.Range("E76").Value = 0
$ws.Range("E77").Value = 2
$ws.Range("E78").Value = 2
$ws.Range("E79").Value = 0

# AttackType
$ws.Range("F65").Value = 'melee'
$ws.Range("F66").Value = 'melee'
$ws.Range("F67").Value = 'melee'
$ws.Range("F68").Value = 'melee'
$ws.Range("F69").Value = 'melee'
$ws.Range("F70").Value = 'ranged'
$ws.Range("F71").Value = 'ranged'
$ws.Range("F73").Value = 'ranged'
$ws.Range("F74").Value = 'melee'
$ws.Range("F75").Value = 'melee'
$ws.Range("F77").Value = 'melee'
$ws.Range("F78").Value = 'melee'

# Speed
$ws.Range("G64").Value = 1
$ws.Range("G65").Value = 4
$ws.Range("G66").Value = 4
$ws.Range("G67").Value = 1
$ws.Range("G68").Value = 1
$ws.Range("G69").Value = 2
$ws.Range("G70").Value = 2
$ws.Range("G71").Value = 1
$ws.Range("G72").Value = 0
$ws.Range("G73").Value = 3
$ws.Range("G74").Value = 2
$ws.Range("G75").Value = 1
$ws.Range("G76").Value = 2
$ws.Range("G77").Value = 3
$ws.Range("G78").Value = 2
$ws.Range("G79").Value = 0

# Health
$ws.Range("H64").Value = 4
$ws.Range("H65").Value = 2
$ws.Range("H66").Value = 5
$ws.Range("H67").Value = 3
$ws.Range("H68").Value = 4
$ws.Range("H69").Value = 5
$ws.Range("H70").Value = 3
$ws.Range("H71").Value = 2
$ws.Range("H72").Value = 0
$ws.Range("H73").Value = 4
$ws.Range("H74").Value = 5
$ws.Range("H75").Value = 7
$ws.Range("H76").Value = 1
$ws.Range("H77").Value = 5
$ws.Range("H78").Value = 9
$ws.Range("H79").Value = 0

# Armor
$ws.Range("I64").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("I66").Value = 1
$ws.Range("I67").Value = 1
$ws.Range("I68").Value = 0
$ws.Range("I69").Value = 1
$ws.Range("I70").Value = 1
$ws.Range("I71").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("I74").Value = 1
$ws.Range("I75").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("I78").Value = 4
$ws.Range("I79").Value = 0

# Ability1
$ws.Range("J64").Value = 'tank heal'
$ws.Range("J65").Value = 'sneak'
$ws.Range("J66").Value = 'inspire'
$ws.Range("J67").Value = 'reach'
$ws.Range("J68").Value = 'opportunity'
$ws.Range("J70").Value = 'snipe'
$ws.Range("J72").Value = 'armor+'
$ws.Range("J74").Value = 'shield'
$ws.Range("J75").Value = 'void'
$ws.Range("J76").Value = 'protect'
$ws.Range("J77").Value = 'reach'
$ws.Range("J78").Value = 'taunt'
$ws.Range("J79").Value = 'health+'

# Ability2
$ws.Range("K77").Value = 'flying'

# --- Death + Dragon splinter cards (rows 80-97) ---

# Card
$ws.Range("A80").Value = 'Animated Corpse'
$ws.Range("A81").Value = 'Haunted Spider'
$ws.Range("A82").Value = 'Skeleton Assassin'
$ws.Range("A83").Value = 'Spineback Wolf'
$ws.Range("A84").Value = 'Maggots'
$ws.Range("A85").Value = 'Cursed Slimeball'
$ws.Range("A86").Value = 'Giant Scorpion'
$ws.Range("A87").Value = 'Undead Badger'
$ws.Range("A88").Value = 'Zintar Mortalis'
$ws.Range("A89").Value = 'Haunted Spirit'
$ws.Range("A90").Value = 'Twisted Jester'
$ws.Range("A91").Value = 'Undead Priest'
$ws.Range("A92").Value = 'Dark Astronomer'
$ws.Range("A93").Value = 'Bone Golem'
$ws.Range("A94").Value = 'Death Elemental'
$ws.Range("A95").Value = 'Contessa L''ament'
$ws.Range("A96").Value = 'Drake of Arnak'
$ws.Range("A97").Value = 'Naga Assassin'

# Role
$ws.Range("B80").Value = 'monster'
$ws.Range("B81").Value = 'monster'
$ws.Range("B82").Value = 'monster'
$ws.Range("B83").Value = 'monster'
$ws.Range("B84").Value = 'monster'
$ws.Range("B85").Value = 'monster'
$ws.Range("B86").Value = 'monster'
$ws.Range("B87").Value = 'monster'
$ws.Range("B88").Value = 'summoner'
$ws.Range("B89").Value = 'monster'
$ws.Range("B90").Value = 'monster'
$ws.Range("B91").Value = 'monster'
$ws.Range("B92").Value = 'monster'
$ws.Range("B93").Value = 'monster'
$ws.Range("B94").Value = 'monster'
$ws.Range("B95").Value = 'summoner'
$ws.Range("B96").Value = 'summoner'
$ws.Range("B97").Value = 'monster'

# Element
$ws.Range("C80").Value = 'death'
$ws.Range("C81").Value = 'death'
$ws.Range("C82").Value = 'death'
$ws.Range("C83").Value = 'death'
$ws.Range("C84").Value = 'death'
$ws.Range("C85").Value = 'death'
$ws.Range("C86").Value = 'death'
$ws.Range("C87").Value = 'death'
$ws.Range("C88").Value = 'death'
$ws.Range("C89").Value = 'death'
$ws.Range("C90").Value = 'death'
$ws.Range("C91").Value = 'death'
$ws.Range("C92").Value = 'death'
$ws.Range("C93").Value = 'death'
$ws.Range("C94").Value = 'death'
$ws.Range("C95").Value = 'death'
$ws.Range("C96").Value = 'dragon'
$ws.Range("C97").Value = 'dragon'

# ManaCost
$ws.Range("D80").Value = 4
$ws.Range("D81").Value = 3
$ws.Range("D82").Value = 3
$ws.Range("D83").Value = 5
$ws.Range("D84").Value = 3
$ws.Range("D85").Value = 1
$ws.Range("D86").Value = 4
$ws.Range("D87").Value = 2
$ws.Range("D88").Value = 3
$ws.Range("D89").Value = 5
$ws.Range("D90").Value = 4
$ws.Range("D91").Value = 2
$ws.Range("D92").Value = 4
$ws.Range("D93").Value = 7
$ws.Range("D94").Value = 3
$ws.Range("D95").Value = 3
$ws.Range("D96").Value = 4
$ws.Range("D97").Value = 2

# Dmg
$ws.Range("E80").Value = 2
$ws.Range("E81").Value = 2
$ws.Range("E82").Value = 1
$ws.Range("E83").Value = 1
$ws.Range("E84").Value = 1
$ws.Range("E85").Value = 1
$ws.Range("E86").Value = 1
$ws.Range("E87").Value = 1
$ws.Range("E88").Value = 0
$ws.Range("E89").Value = 2
$ws.Range("E90").Value = 2
$ws.Range("E91").Value = 0
$ws.Range("E92").Value = 2
$ws.Range("E93").Value = 3
$ws.Range("E94").Value = 1
$ws.Range("E95").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("E97").Value = 1

# AttackType
$ws.Range("F80").Value = 'melee'
$ws.Range("F81").Value = 'ranged'
$ws.Range("F82").Value = 'melee'
$ws.Range("F83").Value = 'melee'
$ws.Range("F84").Value = 'melee'
$ws.Range("F85").Value = 'melee'
$ws.Range("F86").Value = 'melee'
$ws.Range("F87").Value = 'melee'
$ws.Range("F89").Value = 'melee'
$ws.Range("F90").Value = 'ranged'
$ws.Range("F92").Value = 'ranged'
$ws.Range("F93").Value = 'melee'
$ws.Range("F94").Value = 'magic'
$ws.Range("F97").Value = 'ranged'

# Speed
$ws.Range("G80").Value = 1
$ws.Range("G81").Value = 1
$ws.Range("G82").Value = 4
$ws.Range("G83").Value = 6
$ws.Range("G84").Value = 1
$ws.Range("G85").Value = 1
$ws.Range("G86").Value = 2
$ws.Range("G87").Value = 3
$ws.Range("G88").Value = 0
$ws.Range("G89").Value = 2
$ws.Range("G90").Value = 3
$ws.Range("G91").Value = 1
$ws.Range("G92").Value = 1
$ws.Range("G93").Value = 1
$ws.Range("G94").Value = 3
$ws.Range("G95").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("G97").Value = 5

# Health
$ws.Range("H80").Value = 6
$ws.Range("H81").Value = 2
$ws.Range("H82").Value = 2
$ws.Range("H83").Value = 3
$ws.Range("H84").Value = 1
$ws.Range("H85").Value = 1
$ws.Range("H86").Value = 4
$ws.Range("H87").Value = 1
$ws.Range("H88").Value = 0
$ws.Range("H89").Value = 7
$ws.Range("H90").Value = 4
$ws.Range("H91").Value = 3
$ws.Range("H92").Value = 4
$ws.Range("H93").Value = 6
$ws.Range("H94").Value = 2
$ws.Range("H95").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("H97").Value = 2

# Armor
$ws.Range("I80").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("I83").Value = 1
$ws.Range("I84").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("I86").Value = 1
$ws.Range("I87").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("I93").Value = 2
$ws.Range("I94").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("I97").Value = 0

# Ability1
$ws.Range("J82").Value = 'sneak'
$ws.Range("J83").Value = 'reach'
$ws.Range("J84").Value = 'opportunity'
$ws.Range("J87").Value = 'sneak'
$ws.Range("J88").Value = 'melee-'
$ws.Range("J89").Value = 'heal'
$ws.Range("J90").Value = 'snipe'
$ws.Range("J91").Value = 'weaken'
$ws.Range("J93").Value = 'void'
$ws.Range("J94").Value = 'snipe'
$ws.Range("J95").Value = 'ranged-'
$ws.Range("J96").Value = 'armor+'

# Ability2

# Reflect the author's window/view state when they saved: zoomed in
# on the newly-added rows, selection parked on E84.
$excel.ActiveWindow.Zoom = 145
$ws.Range("E84").Select()

